$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell reference -> new value, taken from the refreshed crypto price feed.
$updates = [ordered]@{
    'D2' = '43.468.37'
    'E2' = '  -1.27%  '
    'D3' = '2.374.50'
    'E3' = '  +4.69%  '
    'E4' = '  +0.13%  '
    'D5' = '235.58'
    'E5' = '  +1.08%  '
    'E6' = '  +0.03%  '
    'D7' = '72.32'
    'E7' = '  +13.13%  '
    'E8' = '  +0.07%  '
    'D9' = '0.470'
    'E9' = '  +3.71%  '
    'D10' = '0.0979'
    'E10' = '  -0.32%  '
    'D11' = '56.94'
    'E11' = '  -2.09%  '
    'D12' = '27.29'
    'E12' = '  +1.50%  '
    'D13' = '2.725.12'
    'E13' = '  +4.72%  '
    'D14' = '0.105'
    'E14' = '  -0.78%  '
    'D15' = '15.93'
    'E15' = '  +1.31%  '
    'D16' = '6.32'
    'E16' = '  +2.37%  '
    'D17' = '0.857'
    'D18' = '2.368.28'
    'E18' = '  +4.72%  '
    'D19' = '43.485.60'
    'E19' = '  -1.01%  '
    'D20' = '0.0₂01000'
    'E20' = '  +1.23%  '
    'B21' = 'Litecoin'
    'C21' = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
    'D21' = '74.88'
    'E21' = '  +0.93%  '
    'B22' = 'Uniswap'
    'C22' = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
    'D22' = '6.36'
    'E22' = '  +3.86%  '
    'D23' = '251.56'
    'E23' = '  +0.21%  '
    'D24' = '3.83'
    'E24' = '  +15.40%  '
    'E25' = '  +0.01%  '
    'D26' = '2.48'
    'E26' = '  +1.42%  '
    'B27' = 'EthereumClassic'
    'C27' = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    'D27' = '22.85'
    'E27' = '  +2.69%  '
    'B28' = 'Cosmos'
    'C28' = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    'D28' = '9.99'
    'E28' = '  +0.12%  '
    'B29' = 'Toncoin'
    'C29' = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
    'D29' = '2.22'
    'E29' = '  -3.36%  '
    'D30' = '174.58'
    'E30' = '  +0.46%  '
    'D31' = '1.53'
    'E31' = '  +4.48%  '
    'D32' = '0.129'
    'E32' = '  -7.03%  '
    'E33' = '  +0.14%  '
    'D34' = '5.02'
    'E34' = '  -0.59%  '
    'D35' = '0.0692'
    'D36' = '5.07'
    'E36' = '  +1.25%  '
    'E37' = '  +1.78%  '
    'D38' = '2.45'
    'E38' = '  +6.04%  '
    'E39' = '  -0.75%  '
    'E40' = '  +0.24%  '
    'E41' = '  +0.01%  '
    'B42' = 'FraxShare'
    'C42' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'D42' = '8.89'
    'E42' = '  +1.13%  '
    'B43' = 'InjectiveProtocol'
    'C43' = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
    'D43' = '18.63'
    'E43' = '  +7.10%  '
    'E44' = '  +7.15%  '
    'B45' = 'Aave'
    'C45' = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    'D45' = '100.25'
    'E45' = '  +1.60%  '
    'B46' = 'FTXToken'
    'C46' = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
    'D46' = '4.52'
    'E46' = '  +3.11%  '
    'E47' = '  +1.72%  '
    'D48' = '0.0958'
    'E48' = '  +0.37%  '
    'D49' = '1.443.84'
    'E49' = '  -0.65%  '
    'B50' = 'TerraClassic'
    'C50' = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
    'D50' = '0.000207'
    'E50' = '  -6.92%  '
    'B51' = 'RocketPoolETH'
    'C51' = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
    'D51' = '2.597.73'
    'E51' = '  +4.81%  '
}

# Cells whose new value looks like a plain number (e.g. "235.58") need to be
# forced to stay text -- otherwise Excel would silently reinterpret them as
# numeric cells, which the original sheet never uses for the Price column.
$forceText = @(
    'D5', 'D7', 'D9', 'D10', 'D11', 'D12', 'D14', 'D15', 'D16', 'D17', 'D21', 'D22', 'D23', 'D24', 'D26', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D34', 'D35', 'D36', 'D38', 'D42', 'D43', 'D45', 'D46', 'D48', 'D50'
)

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    if ($forceText -contains $cellRef) {
        $range.Value = "'" + $updates[$cellRef]
        $range.Style = "Normal"
    } else {
        $range.Value = $updates[$cellRef]
    }
}
